# Insert a new data row at row 54 (pushing the existing rows 54-84 down to 55-85)
# and populate the new row with a new "Arveja Verde" price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 54..84 down to 55..85, creating a blank row 54.
$ws.Rows.Item(54).Insert()

# Fill the newly inserted row 54 with the new record.
$ws.Cells.Item(54, 1).Value  = 4
$ws.Cells.Item(54, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(54, 3).Value  = 'Los Lagos'
$ws.Cells.Item(54, 4).Value  = 44518
$ws.Cells.Item(54, 5).Value  = 10
$ws.Cells.Item(54, 6).Value  = 100112022
$ws.Cells.Item(54, 7).Value  = 'Arveja Verde'
$ws.Cells.Item(54, 8).Value  = 'Sin especificar'
$ws.Cells.Item(54, 9).Value  = 'Primera'
$ws.Cells.Item(54, 10).Value = 40
$ws.Cells.Item(54, 11).Value = 20000
$ws.Cells.Item(54, 12).Value = 20000
$ws.Cells.Item(54, 13).Value = 20000
$ws.Cells.Item(54, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(54, 15).Value = 'Región del Maule'
$ws.Cells.Item(54, 16).Value = 800
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = 'Hortaliza'
